# Update the multiplication problems in the practice-sheet table.
# Each original expression is unique within the document, so a simple
# whole-document Find/Replace (MatchWholeWord + MatchCase) targets the
# correct cell every time.
#
# NOTE: "448×3=" appears both as a value to be replaced (from "405×6=")
# and as a value being replaced (to "790×8="). To avoid the newly
# inserted "448×3=" being matched again, the "448×3=" -> "790×8="
# replacement is executed before "405×6=" -> "448×3=".
$d = $word.ActiveDocument

$d.Content.Find.Execute("728×4=", $true, $true, $false, $false, $false, $true, 1, $false, "202×2=", 2) | Out-Null
$d.Content.Find.Execute("779×9=", $true, $true, $false, $false, $false, $true, 1, $false, "930×5=", 2) | Out-Null
$d.Content.Find.Execute("448×3=", $true, $true, $false, $false, $false, $true, 1, $false, "790×8=", 2) | Out-Null  # must run before 405x6 below
$d.Content.Find.Execute("405×6=", $true, $true, $false, $false, $false, $true, 1, $false, "448×3=", 2) | Out-Null
$d.Content.Find.Execute("533×5=", $true, $true, $false, $false, $false, $true, 1, $false, "923×6=", 2) | Out-Null
$d.Content.Find.Execute("365×4=", $true, $true, $false, $false, $false, $true, 1, $false, "449×9=", 2) | Out-Null
$d.Content.Find.Execute("422×4=", $true, $true, $false, $false, $false, $true, 1, $false, "509×8=", 2) | Out-Null
$d.Content.Find.Execute("305×7=", $true, $true, $false, $false, $false, $true, 1, $false, "704×5=", 2) | Out-Null
$d.Content.Find.Execute("531×4=", $true, $true, $false, $false, $false, $true, 1, $false, "874×7=", 2) | Out-Null
$d.Content.Find.Execute("849×6=", $true, $true, $false, $false, $false, $true, 1, $false, "171×9=", 2) | Out-Null
$d.Content.Find.Execute("721×7=", $true, $true, $false, $false, $false, $true, 1, $false, "819×7=", 2) | Out-Null
$d.Content.Find.Execute("837×4=", $true, $true, $false, $false, $false, $true, 1, $false, "151×8=", 2) | Out-Null
$d.Content.Find.Execute("304×5=", $true, $true, $false, $false, $false, $true, 1, $false, "937×2=", 2) | Out-Null
$d.Content.Find.Execute("426×9=", $true, $true, $false, $false, $false, $true, 1, $false, "512×3=", 2) | Out-Null
$d.Content.Find.Execute("896×7=", $true, $true, $false, $false, $false, $true, 1, $false, "622×3=", 2) | Out-Null
$d.Content.Find.Execute("788×3=", $true, $true, $false, $false, $false, $true, 1, $false, "406×5=", 2) | Out-Null
$d.Content.Find.Execute("331×5=", $true, $true, $false, $false, $false, $true, 1, $false, "177×4=", 2) | Out-Null
$d.Content.Find.Execute("914×4=", $true, $true, $false, $false, $false, $true, 1, $false, "224×7=", 2) | Out-Null
$d.Content.Find.Execute("152×8=", $true, $true, $false, $false, $false, $true, 1, $false, "123×7=", 2) | Out-Null
$d.Content.Find.Execute("433×2=", $true, $true, $false, $false, $false, $true, 1, $false, "337×5=", 2) | Out-Null
$d.Content.Find.Execute("460×8=", $true, $true, $false, $false, $false, $true, 1, $false, "537×7=", 2) | Out-Null
$d.Content.Find.Execute("135×5=", $true, $true, $false, $false, $false, $true, 1, $false, "236×3=", 2) | Out-Null
$d.Content.Find.Execute("643×4=", $true, $true, $false, $false, $false, $true, 1, $false, "264×2=", 2) | Out-Null
$d.Content.Find.Execute("567×8=", $true, $true, $false, $false, $false, $true, 1, $false, "999×4=", 2) | Out-Null
$d.Content.Find.Execute("807×9=", $true, $true, $false, $false, $false, $true, 1, $false, "667×9=", 2) | Out-Null
